# Apply the Cover letter edits via Word COM-interop find/replace.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- Main "Deploying swarms..." paragraph -----------------------------
Replace-Text "Deploying swarms" "Deploying teams"

# Merge paragraph breaks within the abstract-like paragraph into spaces.
# (The original runs contain literal newlines inside <w:t>, reflowing the
# text into one continuous paragraph.)
Replace-Text "environment.`nOne key challenge" "environment. One main challenge"
Replace-Text "One main challenge for effective estimation is limited" "One main challenge for effective distributed estimation is the limited"
Replace-Text "resource of robots.`nThis paper presents a distributed estimation scheme, called distributed Bayesian filtering, that is based on" "resource of robots. This paper presents a distributed Bayesian filtering (DBF) scheme that is based on"
Replace-Text "neighboring robots.`nThe proposed scheme can significantly reduce" "neighboring robots. DBF can significantly reduce"
Replace-Text "realistic applications.`nIn addition, theoretical analysis ensures" "realistic applications. In addition, the proposed scheme ensures"
Replace-Text "effectiveness of the proposed method." "effectiveness of the DBF."

# --- Applications sentence ---------------------------------------------
Replace-Text "Applications are 4 and 5, respectively." "Applications are 4 (Safety, Security, Rescue) and 5 (Transportation), respectively."

# --- "Yours sincerely," split around the _GoBack bookmark --------------
# The trailing comma moves to its own run, now located after the
# (empty) _GoBack bookmark instead of before it.
Replace-Text "Yours sincerely," "Yours sincerely"
$t = $d.Content.Text
$idx = $t.IndexOf("Yours sincerely")
$insertPoint = $d.Range($idx + 15, $idx + 15)
$insertPoint.InsertAfter(",") | Out-Null
$newBmRange = $d.Range($idx + 15, $idx + 15)
$d.Bookmarks.Add("_GoBack", $newBmRange) | Out-Null
